$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 15:16"

# --- Refresh per-country COVID-19 stats. The source feed updated a handful of
#     countries totals; the table (rows 4:189) is kept sorted descending by
#     "Casos totales" (col B), so re-applying the sort after the data refresh
#     shifts several rows. Values below are the final, already-resorted rows.
$updates = @{
    9 = @("Estados Unidos", 19781, 398, 147, 19358, 64, 20, 276)
    16 = @("Austria", 2814, 165, 9, 2797, 15, 2, 8)
    17 = @("Noruega", 2003, 44, 1, 1995, 28, 0, 7)
    21 = @("Malasia", 1183, 153, 114, 1061, 26, 5, 8)
    25 = @("Brasil", 987, 17, 2, 973, 18, 1, 12)
    26 = @("Chequia", 925, 92, 5, 920, 7, 0, 0)
    31 = @("Luxemburgo", 670, 186, 6, 656, 3, 3, 8)
    33 = @("Chile", 537, 103, 6, 531, 7, 0, 0)
    34 = @("Finlandia", 521, 71, 10, 510, 2, 1, 1)
    35 = @("Grecia", 495, 0, 19, 463, 20, 3, 13)
    36 = @("Islandia", 473, 64, 5, 468, 1, 0, 0)
    37 = @("Catar", 470, 0, 10, 460, 6, 0, 0)
    38 = @("Polonia", 452, 27, 13, 434, 3, 0, 5)
    39 = @("Indonesia", 450, 81, 20, 392, 0, 6, 38)
    48 = @("Rusia", 306, 53, 16, 289, 0, 0, 1)
    49 = @("Barein", 305, 7, 125, 179, 4, 0, 1)
    50 = @("Egipto", 285, 0, 42, 235, 0, 0, 8)
    51 = @("India", 283, 34, 23, 255, 0, 0, 5)
    52 = @("Hong Kong", 273, 17, 98, 171, 4, 0, 4)
    53 = @("Peru", 263, 0, 1, 258, 5, 0, 4)
    74 = @("Argelia", 95, 1, 43, 40, 0, 1, 12)
    104 = @("Uzbekistan", 41, 8, 0, 41, 0, 0, 0)
    105 = @("Ucrania", 41, 0, 1, 37, 0, 0, 3)
    112 = @("Consejo Danes para los Refugiados", 23, 5, 0, 23, 0, 0, 0)
}

foreach ($r in $updates.Keys) {
    $rowData = $updates[$r]
    $arr = New-Object "object[,]" 1,8
    for ($i = 0; $i -lt 8; $i++) {
        $arr[0,$i] = $rowData[$i]
    }
    $ws.Range("A" + $r + ":H" + $r).Value = $arr
}

Write-Host "Applied" $updates.Count "row updates"
